$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header D1 from "Position" to "PositionSize"
$ws.Range("D1").Value = "PositionSize"

# Widen column D to fit new header text (closest the engine's pixel
# quantization of ColumnWidth can reach to the authored 12.75 chars)
$ws.Columns.Item(4).ColumnWidth = 12

# Add the new trading record for 2021/11/11 (serial 44511)
$ws.Cells.Item(16, 1).Value = 44511
$ws.Cells.Item(16, 1).NumberFormat = $ws.Cells.Item(15, 1).NumberFormat
$ws.Cells.Item(16, 2).Value = 6170
$ws.Cells.Item(16, 3).Value = "short"
$ws.Cells.Item(16, 4).Value = -117
$ws.Cells.Item(16, 5).Value = 47

# Move the active selection to D2, as recorded after the edit
$ws.Range("D2").Select() | Out-Null
